$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, "ВердиоГаст® Растительный комплекс для улучшения пищеварения (БАД ),  капсулы", 74464),
    @(3, "Сб. Фитонефрол (Урологический сбор) 50г", 12449),
    @(4, "Сб. Грудной №4 50г", 21253),
    @(5, "Крушина кора 50г", 7185),
    @(6, "Пустырник трава 50г", 8384),
    @(7, "Чага (березовый гриб) 50г", 15768),
    @(8, "Сенна листья 50г", 16920),
    @(9, "Мать-и-мачеха листья 35г", 22661),
    @(10, "Полынь горькая трава 50г", 30569),
    @(11, "Тысячелистник трава 50г", 11532),
    @(12, "Лен семена 100г", 56011),
    @(13, "Сб. Фитогепатол №2 (Желчегонный сбор №2) 35г", 5041),
    @(14, "Эрва шерстистая трава 30г", 13353),
    @(15, "Ламинарии слоевища (морская капуста) 100г", 14317),
    @(16, "Подорожник большой листья 50г", 9381),
    @(17, "Можжевельник плоды 50г", 10653),
    @(18, "Мята перечная листья 50г", 18228),
    @(19, "Чабрец трава 50г", 19873),
    @(20, "Солодка корни 50г", 32558),
    @(21, "Зверобой трава 50г", 29818),
    @(22, "Пижма цветки 75г", 18135),
    @(23, "Аир корневища 75г", 8007),
    @(24, "Девясил корневища и корни 50г", 18032),
    @(25, "Крапива листья 50г", 15588),
    @(26, "Сб. Фитопектол №2 (Грудной сбор №2) 35г", 9831),
    @(27, "Ромашка цветки вн 50г", 122652),
    @(28, "Валериана корневища с корнями 50г", 21030),
    @(29, "Дуба кора 75г", 78640),
    @(30, "Укроп пахучий плоды 50г", 63029),
    @(31, "Ноготки цветки 50г", 26862),
    @(32, "Бессмертник песчаный цветки 30г", 30084),
    @(33, "Липа цветки 35г", 22859),
    @(34, "Эвкалипт прутовидный листья 75г", 27296),
    @(35, "Багульник болотный побеги 50г", 19611),
    @(36, "Чистотел трава 50г", 23898),
    @(37, "Боярышник плоды 75г", 26265),
    @(38, "Кукуруза столбики с рыльцами 40г", 28385),
    @(39, "Сб. Фитопектол №1 (Грудной сбор №1) 35г", 8670),
    @(40, "Шалфей листья 50г", 46819),
    @(41, "Брусника листья 50г", 24340),
    @(42, "Алтей корни 75г", 8203),
    @(43, "Шиповник плоды низковитаминные 50г", 37808),
    @(44, "Толокнянка листья 50г", 10861),
    @(45, "Череда трава 50г", 21779),
    @(46, "Спорыш трава 50г", 22038),
    @(47, "Рябина плоды 50г", 3150),
    @(48, "Береза почки 50г", 31733),
    @(49, "Фп Детский травяной чай `"ФармаЦветик® для животика`" 20х1,5 г", 1980),
    @(50, "Фп Фиточай `"Лактафитол`" (БАД) 20х1,5 г", 8788),
    @(51, "Фп Детский травяной чай `"ФармаЦветик® для спокойного сна`" 20х1,5 г", 3890),
    @(52, "Фп Детский травяной чай `"ФармаЦветик® для иммунитета`" 20х1,5 г", 3490),
    @(53, "Фп Детский травяной чай `"ФармаЦветик®  при простуде`" 20х1,5 г", 5520),
    @(54, "Фп `"ВердиоГаст® Фиточай для улучшения пищеварения с зеленым чаем`"(БАД) 20*1,5г", 3660),
    @(55, "Фп `"ВердиоГаст® Фиточай для улучшения пищеварения с черным чаем`" (БАД) 20*1,5г", 6670),
    @(56, "Фп Сб. Фитогепатол №3 (Желчегонный сбор №3) 20x2,0г", 30260),
    @(57, "Фп Мята перечная листья 20x1,5г", 23130),
    @(58, "Фп Зверобой трава 20x1,5г", 16259),
    @(59, "Фп Сб. Фитогастрол (Желудочно-кишечный сбор) 20x2,0г", 34326),
    @(60, "Фп `"Щедрость природы® Фиточай успокоительный`"20х2,0 г", 2602),
    @(61, "Фп `"Щедрость природы® Фиточай для иммунитета`" 20х2,0 г", 2718),
    @(62, "Фп Сб. Фитонефрол (Урологический сбор) 20x2,0г", 89728),
    @(63, "Фп `"Щедрость природы® Фиточай кардиологический`" 20х2,0 г", 3870),
    @(64, "Фп Сб. Бруснивер 20x2,0г", 100602),
    @(65, "Фп `"Щедрость природы® Фиточай диабетический`" 20х2,0 г", 3646),
    @(66, "Фп Подорожник листья 20x1,5г", 14780),
    @(67, "Фп Толокнянка листья 20x1,5г", 16730),
    @(68, "Фп Мелисса лекарственная трава 20x1,5г", 19456),
    @(69, "Фп Сб. Грудной №4 20x2,0г", 580081),
    @(70, "Фп Сенна листья 20x1,5г", 41778),
    @(71, "Фп `"Щедрость природы® Фиточай очищающий`" 20х2,0 г", 3392),
    @(72, "Фп Сб. Фитоседан №3 (Успокоительный сбор №3) 20х2,0г", 37527),
    @(73, "Фп Сб. Фитоседан №2 (Успокоительный сбор №2) 20x2,0г", 27251),
    @(74, "Фп Пижма цветки 20х1,5г", 3432),
    @(75, "Фп Сб. Проктофитол (Противогеморроидальный сбор) 20х2,0г", 14634),
    @(76, "Фп Чистотел трава 20х1,5г", 20619),
    @(77, "Фп Хвощ полевой трава 20х1,5г", 18205),
    @(78, "Фп Шалфей листья 20х1,5г", 92142),
    @(79, "Фп Ольха соплодия 20х1,5г", 1724),
    @(80, "Фп Душица трава 20x1,5г", 16614),
    @(81, "Фп Боярышник плоды 20х3,0г", 8498),
    @(82, "Фп Чабрец трава 20x1,5 г", 54104),
    @(83, "Фп Брусника листья 20х1,5г", 57042),
    @(84, "Фп Липа цветки 20x1,5г", 56641),
    @(85, "Фп Сб. Желудочный №3 20x2,0г", 15768),
    @(86, "Фп Сб. Элекасол 20x2,0г", 30384),
    @(87, "Фп Сб. Арфазетин-Э 20x2,0г", 35919),
    @(88, "Фп Пустырник трава 20x1,5г", 27126),
    @(89, "Фп Крапива листья 20x1,5г", 44573),
    @(90, "Фп Пастушья сумка трава 20х1,5г", 4876),
    @(91, "Фп Шиповник плоды 20х2,0г", 36642),
    @(92, "Фп Череда трава 20х1,5г", 44243),
    @(93, "Фп `"Щедрость природы® Фиточай при простуде`" 20х2,0 г", 4266),
    @(94, "Фп `"Щедрость природы® Фиточай для пищеварения`" 20х2,0 г", 2016),
    @(95, "Фп Береза листья 20x1,5г", 4389),
    @(96, "Фп Ромашка цветки 20x1,5г", 1436508),
    @(97, "Фп Золототысячник трава 20х1,5г", 8182),
    @(98, "Фп Фиалка трехцветная трава 20x1,5г", 10044),
    @(99, "Фп Аир корневища 20x1,5г", 12460),
    @(100, "Фп Почечный чай листья 20x1,5г", 48710),
    @(101, "Фп Тысячелистник трава 20x1,5г", 17478),
    @(102, "Фп Кровохлебка корневища и корни 20x1,5г", 5866),
    @(103, "Фп Ноготки цветки 20x1,5г", 41684),
    @(104, "Фп Крушина кора 20x1,5г", 11664),
    @(105, "Фп Валериана корневища с корнями 20x1,5г", 22092),
    @(106, "Фп Дуб кора 20х1,5г", 7245),
    @(107, "Фп Бадан корневища 20x1,5г", 3931),
    @(108, "Фп Лапчатка корневища 20x2,5г", 5310),
    @(109, "Фп Девясил корневища и корни 20х1,5г", 25316)
)

foreach ($item in $data) {
    $r = $item[0]
    $name = $item[1]
    $val = $item[2]
    $ws.Cells.Item($r, 1).Value2 = $name
    $ws.Cells.Item($r, 2).Value2 = $val
}

# Fix the style of the former B74 cell: it used a one-off "0" integer format (numFmtId 1);
# normalize it to the standard "#,##0" thousands format shared by the rest of column B,
# which collapses the now-redundant cellXfs entry.
$ws.Range("B74").NumberFormat = "#,##0"

# Move the active selection from A88 to A90 (topLeftCell / scroll position is unaffected).
$ws.Range("A90").Select()

Write-Output "done"
